$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The author inserted one new row just above the "10 Digit Codes Introduced
# from 1 January 2022" section (originally row 38, now row 39) to record a
# new code. Inserting a whole row at 37 shifts every following row down by
# one, matching the renumbering seen throughout the diff (old row 38 -> new
# row 39, ... old row 70 -> new row 71) and extends the used range from
# A1:B70 to A1:B71.
$ws.Rows("37").Insert(-4121)  # xlShiftDown

# New data for the inserted row.
$ws.Range("A37").Value = "7306 5020 00"
$ws.Range("B37").Value = "Precision Tubes"

# Match the look of the sheet's other "code / description" entries that use
# this same centered, wrapped, small-Arial formatting (e.g. the existing
# "Umbrella Subheading" rows further down the sheet).
$codeCell = $ws.Range("A37")
$codeCell.NumberFormat = "@"
$codeCell.Font.Name = "Arial"
$codeCell.Font.Size = 9
$codeCell.HorizontalAlignment = -4108  # xlCenter
$codeCell.VerticalAlignment = -4108    # xlCenter
$codeCell.WrapText = $true

$descCell = $ws.Range("B37")
$descCell.NumberFormat = "@"
$descCell.Font.Name = "Arial"
$descCell.Font.Size = 9
$descCell.Font.Color = 4473924         # RGB(68,68,68) / FF444444
$descCell.HorizontalAlignment = 1      # xlGeneral
$descCell.VerticalAlignment = -4108    # xlCenter
$descCell.WrapText = $true

# The workbook was left scrolled to the top with B37 selected (no more
# topLeftCell override, and the selection moved from A63 to the newly
# entered cell).
[void]$ws.Range("B37").Select()
